$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet to "Budget" and add a new "Criteria" sheet after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Budget"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Criteria"

# --- Populate the Criteria sheet with budget-name / criteria pairs ---
$ws2.Range("A1").Value = "BUDGET_NAME"
$ws2.Range("B1").Value = "CRITERIA"
$ws2.Range("A2").Value = "Sample Budget 1"
$ws2.Range("B2").Value = "[INTERSTATE]=|Y| AND [INTERNETREPORT]=|State|"
$ws2.Range("A3").Value = "Sample Budget 2"
$ws2.Range("B3").Value = "[INTERSTATE]='Y' AND [INTERNETREPORT]='State'"

# Give the CRITERIA column cells (B2:B3) their own explicit (but still plain)
# cell style, matching the extra cellXfs entry introduced in the workbook.
$ws2.Range("B2:B3").ShrinkToFit = $false

# --- Column widths on the Criteria sheet ---
$ws2.Columns.Item(1).ColumnWidth = 15.85546875
$ws2.Columns.Item(2).ColumnWidth = 47.42578125

# --- Selections on each sheet ---
$ws1.Range("C1").Select()
$ws2.Range("M19").Select()
